$wb = $excel.ActiveWorkbook

# Rename the "Requested quantity" header on the "Weekly Quantity" sheet
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# Rename the "Requested quantity" header on the "Monthly Trend" sheet
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet at the end of the workbook
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# Copy the header formatting (bold, bordered, centered) from the Weekly Quantity header row
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Copy the date-number-format styling from the Weekly Quantity date column
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A61").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header values
$wsForecast.Cells.Item(1,1).Value = "ds"
$wsForecast.Cells.Item(1,2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1,3).Value = "yhat_lower"
$wsForecast.Cells.Item(1,4).Value = "yhat_upper"

# Data rows
$ws = $wsForecast
$ws.Cells.Item(2,1).Value = 45109.99999999999
$ws.Cells.Item(2,2).Value = 122
$ws.Cells.Item(2,3).Value = -835.6488785859967
$ws.Cells.Item(2,4).Value = 1108.635429218274
$ws.Cells.Item(3,1).Value = 45123.99999999999
$ws.Cells.Item(3,2).Value = 146
$ws.Cells.Item(3,3).Value = -775.8332883999547
$ws.Cells.Item(3,4).Value = 1126.534149979286
$ws.Cells.Item(4,1).Value = 45130.99999999999
$ws.Cells.Item(4,2).Value = 158
$ws.Cells.Item(4,3).Value = -807.0533412821774
$ws.Cells.Item(4,4).Value = 1128.336652886524
$ws.Cells.Item(5,1).Value = 45137.99999999999
$ws.Cells.Item(5,2).Value = 170
$ws.Cells.Item(5,3).Value = -766.9727124938161
$ws.Cells.Item(5,4).Value = 1237.583644660089
$ws.Cells.Item(6,1).Value = 45144.99999999999
$ws.Cells.Item(6,2).Value = 182
$ws.Cells.Item(6,3).Value = -851.3232678472796
$ws.Cells.Item(6,4).Value = 1223.868109089246
$ws.Cells.Item(7,1).Value = 45151.99999999999
$ws.Cells.Item(7,2).Value = 194
$ws.Cells.Item(7,3).Value = -714.526498595209
$ws.Cells.Item(7,4).Value = 1157.568101829287
$ws.Cells.Item(8,1).Value = 45158.99999999999
$ws.Cells.Item(8,2).Value = 206
$ws.Cells.Item(8,3).Value = -773.6367651560749
$ws.Cells.Item(8,4).Value = 1220.06883812287
$ws.Cells.Item(9,1).Value = 45165.99999999999
$ws.Cells.Item(9,2).Value = 218
$ws.Cells.Item(9,3).Value = -658.219910387439
$ws.Cells.Item(9,4).Value = 1250.701906606763
$ws.Cells.Item(10,1).Value = 45179.99999999999
$ws.Cells.Item(10,2).Value = 242
$ws.Cells.Item(10,3).Value = -751.9393776315899
$ws.Cells.Item(10,4).Value = 1249.868147522119
$ws.Cells.Item(11,1).Value = 45207.99999999999
$ws.Cells.Item(11,2).Value = 290
$ws.Cells.Item(11,3).Value = -734.833933539636
$ws.Cells.Item(11,4).Value = 1355.762912557114
$ws.Cells.Item(12,1).Value = 45214.99999999999
$ws.Cells.Item(12,2).Value = 302
$ws.Cells.Item(12,3).Value = -665.376520828181
$ws.Cells.Item(12,4).Value = 1301.686758800375
$ws.Cells.Item(13,1).Value = 45221.99999999999
$ws.Cells.Item(13,2).Value = 313
$ws.Cells.Item(13,3).Value = -781.7987796075039
$ws.Cells.Item(13,4).Value = 1207.365516141021
$ws.Cells.Item(14,1).Value = 45228.99999999999
$ws.Cells.Item(14,2).Value = 325
$ws.Cells.Item(14,3).Value = -637.2793939417709
$ws.Cells.Item(14,4).Value = 1390.479132855456
$ws.Cells.Item(15,1).Value = 45249.99999999999
$ws.Cells.Item(15,2).Value = 361
$ws.Cells.Item(15,3).Value = -642.0299437643873
$ws.Cells.Item(15,4).Value = 1405.523230476082
$ws.Cells.Item(16,1).Value = 45270.99999999999
$ws.Cells.Item(16,2).Value = 397
$ws.Cells.Item(16,3).Value = -583.108654658096
$ws.Cells.Item(16,4).Value = 1383.311589361685
$ws.Cells.Item(17,1).Value = 45277.99999999999
$ws.Cells.Item(17,2).Value = 409
$ws.Cells.Item(17,3).Value = -551.4754735673458
$ws.Cells.Item(17,4).Value = 1316.297229367803
$ws.Cells.Item(18,1).Value = 45298.99999999999
$ws.Cells.Item(18,2).Value = 445
$ws.Cells.Item(18,3).Value = -616.9110018414179
$ws.Cells.Item(18,4).Value = 1412.201959894841
$ws.Cells.Item(19,1).Value = 45312.99999999999
$ws.Cells.Item(19,2).Value = 469
$ws.Cells.Item(19,3).Value = -515.6618600909683
$ws.Cells.Item(19,4).Value = 1503.821000276349
$ws.Cells.Item(20,1).Value = 45319.99999999999
$ws.Cells.Item(20,2).Value = 481
$ws.Cells.Item(20,3).Value = -563.9973341163721
$ws.Cells.Item(20,4).Value = 1441.680880064337
$ws.Cells.Item(21,1).Value = 45326.99999999999
$ws.Cells.Item(21,2).Value = 493
$ws.Cells.Item(21,3).Value = -551.2432341865439
$ws.Cells.Item(21,4).Value = 1450.005290694486
$ws.Cells.Item(22,1).Value = 45333.99999999999
$ws.Cells.Item(22,2).Value = 505
$ws.Cells.Item(22,3).Value = -439.1341906831012
$ws.Cells.Item(22,4).Value = 1488.758360544375
$ws.Cells.Item(23,1).Value = 45340.99999999999
$ws.Cells.Item(23,2).Value = 517
$ws.Cells.Item(23,3).Value = -507.7120782371933
$ws.Cells.Item(23,4).Value = 1503.982739685488
$ws.Cells.Item(24,1).Value = 45347.99999999999
$ws.Cells.Item(24,2).Value = 529
$ws.Cells.Item(24,3).Value = -456.8184637166165
$ws.Cells.Item(24,4).Value = 1567.414708711937
$ws.Cells.Item(25,1).Value = 45354.99999999999
$ws.Cells.Item(25,2).Value = 541
$ws.Cells.Item(25,3).Value = -461.554446824467
$ws.Cells.Item(25,4).Value = 1567.856986887738
$ws.Cells.Item(26,1).Value = 45361.99999999999
$ws.Cells.Item(26,2).Value = 553
$ws.Cells.Item(26,3).Value = -537.3575275739037
$ws.Cells.Item(26,4).Value = 1555.852130177519
$ws.Cells.Item(27,1).Value = 45368.99999999999
$ws.Cells.Item(27,2).Value = 565
$ws.Cells.Item(27,3).Value = -366.6105631300451
$ws.Cells.Item(27,4).Value = 1500.515127638059
$ws.Cells.Item(28,1).Value = 45375.99999999999
$ws.Cells.Item(28,2).Value = 577
$ws.Cells.Item(28,3).Value = -472.4592378694835
$ws.Cells.Item(28,4).Value = 1591.573082093481
$ws.Cells.Item(29,1).Value = 45382.99999999999
$ws.Cells.Item(29,2).Value = 589
$ws.Cells.Item(29,3).Value = -324.7824249712058
$ws.Cells.Item(29,4).Value = 1567.445134144729
$ws.Cells.Item(30,1).Value = 45389.99999999999
$ws.Cells.Item(30,2).Value = 600
$ws.Cells.Item(30,3).Value = -347.3442505267897
$ws.Cells.Item(30,4).Value = 1603.488369598753
$ws.Cells.Item(31,1).Value = 45396.99999999999
$ws.Cells.Item(31,2).Value = 612
$ws.Cells.Item(31,3).Value = -343.4192898161252
$ws.Cells.Item(31,4).Value = 1678.69962025682
$ws.Cells.Item(32,1).Value = 45410.99999999999
$ws.Cells.Item(32,2).Value = 636
$ws.Cells.Item(32,3).Value = -355.7514339507796
$ws.Cells.Item(32,4).Value = 1571.535792387506
$ws.Cells.Item(33,1).Value = 45417.99999999999
$ws.Cells.Item(33,2).Value = 648
$ws.Cells.Item(33,3).Value = -291.4279303591582
$ws.Cells.Item(33,4).Value = 1627.493464240261
$ws.Cells.Item(34,1).Value = 45424.99999999999
$ws.Cells.Item(34,2).Value = 660
$ws.Cells.Item(34,3).Value = -335.2011969468355
$ws.Cells.Item(34,4).Value = 1657.882338788664
$ws.Cells.Item(35,1).Value = 45431.99999999999
$ws.Cells.Item(35,2).Value = 672
$ws.Cells.Item(35,3).Value = -277.7044518082385
$ws.Cells.Item(35,4).Value = 1689.058487734388
$ws.Cells.Item(36,1).Value = 45438.99999999999
$ws.Cells.Item(36,2).Value = 684
$ws.Cells.Item(36,3).Value = -309.1335630895834
$ws.Cells.Item(36,4).Value = 1646.193262876515
$ws.Cells.Item(37,1).Value = 45445.99999999999
$ws.Cells.Item(37,2).Value = 696
$ws.Cells.Item(37,3).Value = -313.5480901082159
$ws.Cells.Item(37,4).Value = 1685.470336564062
$ws.Cells.Item(38,1).Value = 45459.99999999999
$ws.Cells.Item(38,2).Value = 720
$ws.Cells.Item(38,3).Value = -311.9054931381003
$ws.Cells.Item(38,4).Value = 1747.910078258309
$ws.Cells.Item(39,1).Value = 45466.99999999999
$ws.Cells.Item(39,2).Value = 732
$ws.Cells.Item(39,3).Value = -237.3658952688732
$ws.Cells.Item(39,4).Value = 1719.249941882912
$ws.Cells.Item(40,1).Value = 45473.99999999999
$ws.Cells.Item(40,2).Value = 744
$ws.Cells.Item(40,3).Value = -305.2526771235004
$ws.Cells.Item(40,4).Value = 1671.15582977682
$ws.Cells.Item(41,1).Value = 45480.99999999999
$ws.Cells.Item(41,2).Value = 756
$ws.Cells.Item(41,3).Value = -270.1574696990393
$ws.Cells.Item(41,4).Value = 1713.127012938206
$ws.Cells.Item(42,1).Value = 45487.99999999999
$ws.Cells.Item(42,2).Value = 768
$ws.Cells.Item(42,3).Value = -267.5646693625378
$ws.Cells.Item(42,4).Value = 1734.957999439636
$ws.Cells.Item(43,1).Value = 45529.99999999999
$ws.Cells.Item(43,2).Value = 840
$ws.Cells.Item(43,3).Value = -87.135755479343
$ws.Cells.Item(43,4).Value = 1820.451143798118
$ws.Cells.Item(44,1).Value = 45536.99999999999
$ws.Cells.Item(44,2).Value = 852
$ws.Cells.Item(44,3).Value = -123.7996049787787
$ws.Cells.Item(44,4).Value = 1861.673564260666
$ws.Cells.Item(45,1).Value = 45543.99999999999
$ws.Cells.Item(45,2).Value = 864
$ws.Cells.Item(45,3).Value = -50.49763766521973
$ws.Cells.Item(45,4).Value = 1882.277799950448
$ws.Cells.Item(46,1).Value = 45550.99999999999
$ws.Cells.Item(46,2).Value = 876
$ws.Cells.Item(46,3).Value = -81.41000890748454
$ws.Cells.Item(46,4).Value = 1910.218518262751
$ws.Cells.Item(47,1).Value = 45557.99999999999
$ws.Cells.Item(47,2).Value = 887
$ws.Cells.Item(47,3).Value = -63.30010318795716
$ws.Cells.Item(47,4).Value = 1842.390408637348
$ws.Cells.Item(48,1).Value = 45564.99999999999
$ws.Cells.Item(48,2).Value = 899
$ws.Cells.Item(48,3).Value = -70.93364991235931
$ws.Cells.Item(48,4).Value = 1917.985440446598
$ws.Cells.Item(49,1).Value = 45571.99999999999
$ws.Cells.Item(49,2).Value = 911
$ws.Cells.Item(49,3).Value = -31.35824692306864
$ws.Cells.Item(49,4).Value = 1932.104410621156
$ws.Cells.Item(50,1).Value = 45592.99999999999
$ws.Cells.Item(50,2).Value = 947
$ws.Cells.Item(50,3).Value = -61.8802172323444
$ws.Cells.Item(50,4).Value = 1896.308982921634
$ws.Cells.Item(51,1).Value = 45606.99999999999
$ws.Cells.Item(51,2).Value = 971
$ws.Cells.Item(51,3).Value = -56.43208876090809
$ws.Cells.Item(51,4).Value = 1943.667231789565
$ws.Cells.Item(52,1).Value = 45613.99999999999
$ws.Cells.Item(52,2).Value = 983
$ws.Cells.Item(52,3).Value = -75.96522114730642
$ws.Cells.Item(52,4).Value = 1921.932301637178
$ws.Cells.Item(53,1).Value = 45627.99999999999
$ws.Cells.Item(53,2).Value = 1007
$ws.Cells.Item(53,3).Value = 107.7730327275476
$ws.Cells.Item(53,4).Value = 2047.107929252176
$ws.Cells.Item(54,1).Value = 45634.99999999999
$ws.Cells.Item(54,2).Value = 1019
$ws.Cells.Item(54,3).Value = -11.65300201435021
$ws.Cells.Item(54,4).Value = 2000.963471104664
$ws.Cells.Item(55,1).Value = 45641.99999999999
$ws.Cells.Item(55,2).Value = 1031
$ws.Cells.Item(55,3).Value = 19.46228475037651
$ws.Cells.Item(55,4).Value = 2086.701711839247
$ws.Cells.Item(56,1).Value = 45648.99999999999
$ws.Cells.Item(56,2).Value = 1043
$ws.Cells.Item(56,3).Value = 114.0327951148185
$ws.Cells.Item(56,4).Value = 2008.496174132569
$ws.Cells.Item(57,1).Value = 45655.99999999999
$ws.Cells.Item(57,2).Value = 1055
$ws.Cells.Item(57,3).Value = 112.6974689427577
$ws.Cells.Item(57,4).Value = 2133.60344434086
$ws.Cells.Item(58,1).Value = 45662.99999999999
$ws.Cells.Item(58,2).Value = 1067
$ws.Cells.Item(58,3).Value = 29.06677369992552
$ws.Cells.Item(58,4).Value = 2072.228942539828
$ws.Cells.Item(59,1).Value = 45669.99999999999
$ws.Cells.Item(59,2).Value = 1079
$ws.Cells.Item(59,3).Value = 72.29990649589411
$ws.Cells.Item(59,4).Value = 2089.928648386386
$ws.Cells.Item(60,1).Value = 45676.99999999999
$ws.Cells.Item(60,2).Value = 1091
$ws.Cells.Item(60,3).Value = 81.04832189973969
$ws.Cells.Item(60,4).Value = 2076.255375539514
$ws.Cells.Item(61,1).Value = 45683.99999999999
$ws.Cells.Item(61,2).Value = 1103
$ws.Cells.Item(61,3).Value = 103.0072665117113
$ws.Cells.Item(61,4).Value = 2072.705322274054

